$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update header row of sheet1 (o_10): add evaluator_partial_correctness column ---
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("D1").Copy() | Out-Null
$ws1.Range("E1").PasteSpecial(-4122) | Out-Null

# --- Update data row of sheet1 (o_10) ---
$text5 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node J? Return the sequence of nodes in response.
   A B C D E F G H I J
 A 0 1 1 0 0 0 0 0 0 0
 B 1 0 0 0 0 0 0 0 0 0
 C 1 0 0 1 1 0 0 0 0 0
 D 0 0 1 0 0 0 0 0 0 0
 E 0 0 1 0 0 1 0 0 0 0
 F 0 0 0 0 1 0 1 0 1 0
 G 0 0 0 0 0 1 0 1 0 0
 H 0 0 0 0 0 0 1 0 0 0
 I 0 0 0 0 0 1 0 0 0 1
 J 0 0 0 0 0 0 0 0 1 0
    
"@
$ws1.Range("A2").Value = $text5
$ws1.Range("B2").Value = "A -> C -> E -> F -> I -> J"
$ws1.Range("C2").Value = "The shortest path from node A to node J is A -> C -> E -> F -> I -> J."
$ws1.Range("D2").Value = "invalid input"
$ws1.Range("E2").Value = "6/6"
$ws1.Rows.Item(2).AutoFit() | Out-Null

# --- Add sheet o_20 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"
$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:E1").Copy() | Out-Null
$ws2.Range("A1:E1").PasteSpecial(-4122) | Out-Null

$text10 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node T? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 1 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 B 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 J 1 0 0 0 0 0 0 0 0 0 1 1 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 1 1 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
"@
$ws2.Range("A2").Value = $text10
$ws2.Range("B2").Value = "A -> J -> L -> M -> Q -> R -> S -> T"
$ws2.Range("C2").Value = "The shortest path from node A to node T is: A -> C -> D -> E -> F -> H -> I -> R -> S -> T."
$ws2.Range("D2").Value = "invalid input"
$ws2.Range("E2").Value = "1/8"
$ws2.Rows.Item(2).AutoFit() | Out-Null

# --- Add sheet o_20_jumbled ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"
$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:E1").Copy() | Out-Null
$ws3.Range("A1:E1").PasteSpecial(-4122) | Out-Null

$text14 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node T? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 1 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
"@
$ws3.Range("A2").Value = $text14
$ws3.Range("B2").Value = "A -> B -> G -> H -> K -> L -> Q -> R -> S -> T"
$ws3.Range("C2").Value = "The shortest path from node A to node T is: A - B - G - H - K - L - Q - R - S - T."
$ws3.Range("D2").Value = "invalid input"
$ws3.Range("E2").Value = "10/10"
$ws3.Rows.Item(2).AutoFit() | Out-Null

# --- Restore active sheet to o_10 ---
$ws1.Activate()

Write-Output "done"
